# Generate Report for Handoff
# Updates the handoff timestamps for the "ac6590b2-9e2b-48a3-9eff-45886edc2c3b" file
# (row 6 on every sheet) to reflect a new handoff/xliff-generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-22 22:41:35"

# --- zh-cn sheet ------------------------------------------------------
# Column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-22 22:41:30"

# --- de-de sheet ------------------------------------------------------
# Column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-22 22:41:35"
